$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a purely-numeric-looking string into a cell as genuine
# TEXT, without Excel's usual "looks like a number -> store as number"
# auto-coercion, and *without* disturbing the cell's existing style
# (a plain quote-prefixed/Text-formatted assignment would mark the
# style with quotePrefix/custom numFmt and allocate a brand new cellXf).
# Trick: compute the literal digits via TEXT() in a scratch cell, copy
# it, then PasteSpecial *values only* (xlPasteValues = -4163) onto the
# destination - this carries over the String type but leaves the
# destination cell's own formatting/style untouched.
function Set-PlainTextNumber($ws, $addr, $digits) {
    $scratch = $ws.Range("Z1")
    $scratch.Formula = "=TEXT(" + $digits + ",""0"")"
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.Clear()
}

# --- Row 2 ---
$ws.Range("B2").Value = "тск"
$ws.Range("C2").Value = "вспашка"
$ws.Range("D2").Value = "кукурузу"
Set-PlainTextNumber $ws "E2" "70"
$ws.Range("F2").Value = "Не указано"
$ws.Range("G2").Value = "Не указано"
$ws.Range("H2").Value = "Не указано"

# --- Row 3 ---
$ws.Range("B3").Value = "Не указано"
$ws.Range("C3").Value = "100 выравнивание зяби"
$ws.Range("D3").Value = "сою"
Set-PlainTextNumber $ws "E3" "155"
Set-PlainTextNumber $ws "F3" "1377"
$ws.Range("G3").Value = "Не указано"
$ws.Range("H3").Value = "Не указано"

$excel.CutCopyMode = $false
